# Auto-generated Excel COM-interop script to update cryptos list
# Commit: "Updated cryptos list on Wed Sep 18 16:19:22 UTC 2024 with GitHub Actions"
# Source data are text cells (inline strings in the original OOXML); many of the
# updated Price values look numeric (e.g. "539.82"), so plain assignment would make
# Excel auto-convert them to real numbers. To preserve them as text we temporarily
# force a Text number format, assign the value, then restore the default style so
# no visible formatting change remains (matches the original unstyled cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.586.21"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").Value = "2.301.89"
$ws.Range("E3").Value = "  -3.41%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.58%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("D9").Value = "2.297.73"
$ws.Range("E9").Value = "  -3.50%  "
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.150"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").Value = "2.711.23"
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.61%  "
$ws.Range("D16").Value = "59.450.21"
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.01%  "
$ws.Range("D18").Value = "2.294.63"
$ws.Range("E18").Value = "  -3.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "308.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.30%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  -3.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").Value = "0.0₃0711"
$ws.Range("E32").Value = "  -5.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("B34").Value = "PolygonEcosystemToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.375"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -7.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.97%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "309.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("E42").Value = "  -5.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0933"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.564"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0487"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.22%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0211"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0215"
$ws.Range("E50").Value = "  +6.03%  "
$ws.Range("E51").Value = "  -0.51%  "
